$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 60, shifting existing rows 60-62 down to 61-63.
$ws.Rows.Item(60).Insert()

# Copy the style used by column D (date column) down into the new row.
$ws.Cells.Item(60, 4).Value = $ws.Cells.Item(61, 4).Value

# Fill the new row 60 with the "Navel Late" / "Tercera" record.
$ws.Cells.Item(60, 1).Value = 1
$ws.Cells.Item(60, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(60, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(60, 4).Value = 44448
$ws.Cells.Item(60, 5).Value = 15
$ws.Cells.Item(60, 6).Value = "Fruta"
$ws.Cells.Item(60, 7).Value = 100102
$ws.Cells.Item(60, 8).Value = "Cítricos"
$ws.Cells.Item(60, 9).Value = 100102005
$ws.Cells.Item(60, 10).Value = "Naranja"
$ws.Cells.Item(60, 11).Value = "Navel Late"
$ws.Cells.Item(60, 12).Value = "Tercera"
$ws.Cells.Item(60, 13).Value = 200
$ws.Cells.Item(60, 14).Value = 600
$ws.Cells.Item(60, 15).Value = 650
$ws.Cells.Item(60, 16).Value = 625
$ws.Cells.Item(60, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(60, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(60, 19).Value = 625
$ws.Cells.Item(60, 20).Value = 1

# Update the cells that differ between the old (pre-insert) row content and
# the target content, now that rows have shifted down by one.
# Row 61 (was old row 60): date changes from 44448 placeholder to 44167, region to Coquimbo (already correct from old row60 data after shift)
$ws.Cells.Item(61, 4).Value = 44167
$ws.Cells.Item(61, 18).Value = "Región de Coquimbo"
